$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1149.0308
$ws.Range("I15").Value = 1149.0308
$ws.Range("K15").Value = 3447.0924
$ws.Range("M15").Value = -3278.0924
$ws.Range("H98").Value = 1042.7333
$ws.Range("I98").Value = 1085.8334
$ws.Range("J98").Value = 870.3333
$ws.Range("K98").Value = 1085.8334
$ws.Range("L98").Value = 870.3333
$ws.Range("M98").Value = 412.1666
$ws.Range("N98").Value = -3866.3333
$ws.Range("H122").Value = 1042.7333
$ws.Range("I122").Value = 1085.8334
$ws.Range("J122").Value = 870.3333
$ws.Range("K122").Value = 3257.5002
$ws.Range("L122").Value = 2610.9999
$ws.Range("M122").Value = -807.5001999999999
$ws.Range("N122").Value = -7510.9999
$ws.Range("H124").Value = 46500
$ws.Range("J124").Value = 46500
$ws.Range("L124").Value = 46500
$ws.Range("N124").Value = -56320
$ws.Range("H125").Value = 3170.1177
$ws.Range("I125").Value = 2669.3333
$ws.Range("J125").Value = 4372
$ws.Range("K125").Value = 24023.9997
$ws.Range("L125").Value = 39348
$ws.Range("M125").Value = -21563.9997
$ws.Range("N125").Value = -44268
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 44306.168
$ws.Range("I45").Value = 73648.21000000001
$ws.Range("K45").Value = 73648.21000000001
$ws.Range("M45").Value = -73271.21000000001
$ws.Range("H46").Value = 3792
$ws.Range("J46").Value = 4813
$ws.Range("L46").Value = 4813
$ws.Range("N46").Value = -5451
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 97038.336
$ws.Range("I31").Value = 100854.1
$ws.Range("J31").Value = 94793.766
$ws.Range("K31").Value = 100854.1
$ws.Range("L31").Value = 94793.766
$ws.Range("M31").Value = -100559.1
$ws.Range("N31").Value = -95383.766
$ws.Range("H34").Value = 97038.336
$ws.Range("I34").Value = 100854.1
$ws.Range("J34").Value = 94793.766
$ws.Range("K34").Value = 100854.1
$ws.Range("L34").Value = 94793.766
$ws.Range("M34").Value = -100652.1
$ws.Range("N34").Value = -95197.766
$ws.Range("H86").Value = 1972.9231
$ws.Range("I86").Value = 1505.2941
$ws.Range("J86").Value = 2856.2222
$ws.Range("K86").Value = 1505.2941
$ws.Range("L86").Value = 2856.2222
$ws.Range("M86").Value = -382.2941000000001
$ws.Range("N86").Value = -5102.2222
$ws.Range("H89").Value = 1972.9231
$ws.Range("I89").Value = 1505.2941
$ws.Range("J89").Value = 2856.2222
$ws.Range("K89").Value = 7526.4705
$ws.Range("L89").Value = 14281.111
$ws.Range("M89").Value = -1910.4705
$ws.Range("N89").Value = -25513.111
$ws.Range("H94").Value = 1387.3
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1387.3
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1387.3
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2289.3
$ws.Range("H99").Value = 16854
$ws.Range("I99").Value = 2592.8
$ws.Range("J99").Value = 52507
$ws.Range("K99").Value = 2592.8
$ws.Range("L99").Value = 52507
$ws.Range("M99").Value = -1094.8
$ws.Range("N99").Value = -55503
$ws.Range("H126").Value = 16854
$ws.Range("I126").Value = 2592.8
$ws.Range("J126").Value = 52507
$ws.Range("K126").Value = 7778.400000000001
$ws.Range("L126").Value = 157521
$ws.Range("M126").Value = -5308.400000000001
$ws.Range("N126").Value = -162461
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 892.25
$ws.Range("I103").Value = 381.8
$ws.Range("J103").Value = 1256.8572
$ws.Range("K103").Value = 1145.4
$ws.Range("L103").Value = 3770.5716
$ws.Range("M103").Value = -266.4000000000001
$ws.Range("N103").Value = -5528.571599999999
$ws.Range("H105").Value = 7960.3335
$ws.Range("I105").Value = 4500
$ws.Range("J105").Value = 8392.875
$ws.Range("K105").Value = 13500
$ws.Range("L105").Value = 25178.625
$ws.Range("M105").Value = -10879
$ws.Range("N105").Value = -30420.625
$ws.Range("H131").Value = 867.61
$ws.Range("I131").Value = 730
$ws.Range("J131").Value = 869
$ws.Range("K131").Value = 2190
$ws.Range("L131").Value = 2607
$ws.Range("M131").Value = 2850
$ws.Range("N131").Value = -12687
$ws.Range("H134").Value = 2563.725
$ws.Range("I134").Value = 2215.2354
$ws.Range("J134").Value = 2821.3044
$ws.Range("K134").Value = 6645.706200000001
$ws.Range("L134").Value = 8463.913199999999
$ws.Range("M134").Value = -1575.706200000001
$ws.Range("N134").Value = -18603.9132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3817.5
$ws.Range("I102").Value = 3772
$ws.Range("J102").Value = 3893.3333
$ws.Range("K102").Value = 3772
$ws.Range("L102").Value = 3893.3333
$ws.Range("M102").Value = -2150
$ws.Range("N102").Value = -7137.3333
$ws.Range("H107").Value = 505394
$ws.Range("I107").Value = 341.17648
$ws.Range("J107").Value = 3367360
$ws.Range("K107").Value = 341.17648
$ws.Range("L107").Value = 3367360
$ws.Range("M107").Value = 1578.82352
$ws.Range("N107").Value = -3371200
$ws.Range("H122").Value = 951.7778
$ws.Range("I122").Value = 966.7143
$ws.Range("J122").Value = 899.5
$ws.Range("K122").Value = 2900.1429
$ws.Range("L122").Value = 2698.5
$ws.Range("M122").Value = -450.1428999999998
$ws.Range("N122").Value = -7598.5
$ws.Range("H126").Value = 3089.6553
$ws.Range("I126").Value = 2640.9473
$ws.Range("J126").Value = 3942.2
$ws.Range("K126").Value = 7922.841899999999
$ws.Range("L126").Value = 11826.6
$ws.Range("M126").Value = -5452.841899999999
$ws.Range("N126").Value = -16766.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2773.842
$ws.Range("I7").Value = 1800.2307
$ws.Range("J7").Value = 4883.3335
$ws.Range("K7").Value = 1800.2307
$ws.Range("L7").Value = 4883.3335
$ws.Range("M7").Value = -1688.2307
$ws.Range("N7").Value = -5107.3335
$ws.Range("H40").Value = 68693.2
$ws.Range("I40").Value = 334266.34
$ws.Range("K40").Value = 334266.34
$ws.Range("M40").Value = -334130.34
$ws.Range("H126").Value = 2773.842
$ws.Range("I126").Value = 1800.2307
$ws.Range("J126").Value = 4883.3335
$ws.Range("K126").Value = 5400.6921
$ws.Range("L126").Value = 14650.0005
$ws.Range("M126").Value = -2930.6921
$ws.Range("N126").Value = -19590.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10799.667
$ws.Range("J45").Value = 10799.667
$ws.Range("L45").Value = 10799.667
$ws.Range("N45").Value = -11781.667
$ws.Range("H122").Value = 2000.5555
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2801
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 8403
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -13303
$ws.Range("H126").Value = 1729.1538
$ws.Range("I126").Value = 1708.2
$ws.Range("J126").Value = 1799
$ws.Range("K126").Value = 5124.6
$ws.Range("L126").Value = 5397
$ws.Range("M126").Value = -2654.6
$ws.Range("N126").Value = -10337
$ws.Range("H132").Value = 1961.3823
$ws.Range("I132").Value = 1589.8334
$ws.Range("J132").Value = 3394.5
$ws.Range("K132").Value = 4769.5002
$ws.Range("L132").Value = 10183.5
$ws.Range("M132").Value = -2239.5002
$ws.Range("N132").Value = -15243.5
